$d = $word.ActiveDocument

$pairs = @(
    @("33×47=1551", "95×94=8930"),
    @("13×39=507", "83×43=3569"),
    @("42×50=2100", "11×50=550"),
    @("87×51=4437", "83×22=1826"),
    @("59×79=4661", "97×85=8245"),
    @("23×29=667", "52×73=3796"),
    @("34×24=816", "75×92=6900"),
    @("13×42=546", "12×97=1164"),
    @("55×65=3575", "11×90=990"),
    @("92×11=1012", "29×68=1972"),
    @("47×14=658", "97×94=9118"),
    @("83×55=4565", "52×23=1196"),
    @("19×28=532", "29×83=2407"),
    @("16×84=1344", "21×14=294"),
    @("95×57=5415", "65×98=6370"),
    @("95×35=3325", "13×11=143"),
    @("14×57=798", "39×61=2379"),
    @("93×22=2046", "50×32=1600"),
    @("14×96=1344", "53×31=1643"),
    @("39×37=1443", "68×14=952"),
    @("34×12=408", "22×85=1870"),
    @("23×40=920", "97×76=7372"),
    @("67×30=2010", "87×39=3393"),
    @("28×22=616", "82×46=3772"),
    @("65×49=3185", "66×81=5346")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
